# Auto-generated edit script applying the Chocobo_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1210.1428
$ws.Cells.Item(6, 9).Value = 174.2
$ws.Cells.Item(6, 11).Value = 522.5999999999999
$ws.Cells.Item(6, 13).Value = -410.5999999999999
$ws.Cells.Item(15, 8).Value = 182.06
$ws.Cells.Item(15, 9).Value = 182.06
$ws.Cells.Item(15, 11).Value = 546.1800000000001
$ws.Cells.Item(15, 13).Value = -377.1800000000001
$ws.Cells.Item(17, 8).Value = 1643.3433
$ws.Cells.Item(17, 10).Value = 1729.8723
$ws.Cells.Item(17, 12).Value = 5189.6169
$ws.Cells.Item(17, 14).Value = -5525.6169
$ws.Cells.Item(21, 8).Value = 65171.2
$ws.Cells.Item(21, 9).Value = 80019
$ws.Cells.Item(21, 10).Value = 42899.5
$ws.Cells.Item(21, 11).Value = 80019
$ws.Cells.Item(21, 12).Value = 42899.5
$ws.Cells.Item(21, 13).Value = -79551
$ws.Cells.Item(21, 14).Value = -43835.5
$ws.Cells.Item(23, 8).Value = 65171.2
$ws.Cells.Item(23, 9).Value = 80019
$ws.Cells.Item(23, 10).Value = 42899.5
$ws.Cells.Item(23, 11).Value = 80019
$ws.Cells.Item(23, 12).Value = 42899.5
$ws.Cells.Item(23, 13).Value = -79785
$ws.Cells.Item(23, 14).Value = -43367.5
$ws.Cells.Item(28, 8).Value = 540.25
$ws.Cells.Item(28, 9).Value = 519.63635
$ws.Cells.Item(28, 10).Value = 615.8333
$ws.Cells.Item(28, 11).Value = 519.63635
$ws.Cells.Item(28, 12).Value = 615.8333
$ws.Cells.Item(28, 13).Value = -34.63634999999999
$ws.Cells.Item(28, 14).Value = -1585.8333
$ws.Cells.Item(29, 8).Value = 825.25
$ws.Cells.Item(29, 9).Value = 825.25
$ws.Cells.Item(29, 11).Value = 2475.75
$ws.Cells.Item(29, 13).Value = -2194.75
$ws.Cells.Item(38, 8).Value = 2676.682
$ws.Cells.Item(38, 9).Value = 118.7
$ws.Cells.Item(38, 10).Value = 4808.3335
$ws.Cells.Item(38, 11).Value = 356.1
$ws.Cells.Item(38, 12).Value = 14425.0005
$ws.Cells.Item(38, 13).Value = 15.89999999999998
$ws.Cells.Item(38, 14).Value = -15169.0005
$ws.Cells.Item(43, 8).Value = 2440.5
$ws.Cells.Item(43, 9).Value = 1512
$ws.Cells.Item(43, 10).Value = 2797.6155
$ws.Cells.Item(43, 11).Value = 1512
$ws.Cells.Item(43, 12).Value = 2797.6155
$ws.Cells.Item(43, 13).Value = -1443
$ws.Cells.Item(43, 14).Value = -2935.6155
$ws.Cells.Item(58, 8).Value = 10271.462
$ws.Cells.Item(58, 9).Value = 1503.2222
$ws.Cells.Item(58, 10).Value = 30000
$ws.Cells.Item(58, 11).Value = 4509.6666
$ws.Cells.Item(58, 12).Value = 90000
$ws.Cells.Item(58, 13).Value = -4359.6666
$ws.Cells.Item(58, 14).Value = -90300
$ws.Cells.Item(96, 8).Value = 594.2857
$ws.Cells.Item(96, 9).Value = 261.22223
$ws.Cells.Item(96, 10).Value = 1193.8
$ws.Cells.Item(96, 11).Value = 783.66669
$ws.Cells.Item(96, 12).Value = 3581.4
$ws.Cells.Item(96, 13).Value = 589.33331
$ws.Cells.Item(96, 14).Value = -6327.4
$ws.Cells.Item(100, 8).Value = 22223618
$ws.Cells.Item(100, 9).Value = 22223618
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 22223618
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = -22223077
$ws.Cells.Item(100, 14).ClearContents()
$ws.Cells.Item(111, 8).Value = 1326.4
$ws.Cells.Item(111, 10).Value = 1326.4
$ws.Cells.Item(111, 12).Value = 3979.2
$ws.Cells.Item(111, 14).Value = -10113.2
$ws.Cells.Item(112, 8).Value = 11495855
$ws.Cells.Item(112, 9).Value = 333333860
$ws.Cells.Item(112, 10).Value = 1640.2024
$ws.Cells.Item(112, 11).Value = 1000001580
$ws.Cells.Item(112, 12).Value = 4920.607199999999
$ws.Cells.Item(112, 13).Value = -1000000472
$ws.Cells.Item(112, 14).Value = -7136.607199999999
$ws.Cells.Item(115, 8).Value = 1566.125
$ws.Cells.Item(115, 9).Value = 1566.125
$ws.Cells.Item(115, 11).Value = 4698.375
$ws.Cells.Item(115, 13).Value = -3131.375
$ws.Cells.Item(118, 8).Value = 1717.7368
$ws.Cells.Item(118, 9).Value = 1796.6666
$ws.Cells.Item(118, 10).Value = 1702.9375
$ws.Cells.Item(118, 11).Value = 5389.9998
$ws.Cells.Item(118, 12).Value = 5108.8125
$ws.Cells.Item(118, 13).Value = -3732.9998
$ws.Cells.Item(118, 14).Value = -8422.8125
$ws.Cells.Item(129, 8).Value = 905.6667
$ws.Cells.Item(129, 10).Value = 1158.5
$ws.Cells.Item(129, 12).Value = 3475.5
$ws.Cells.Item(129, 14).Value = -13475.5
$ws.Cells.Item(132, 8).Value = 129059.28
$ws.Cells.Item(132, 9).Value = 184359.4
$ws.Cells.Item(132, 11).Value = 553078.2
$ws.Cells.Item(132, 13).Value = -550548.2
$ws.Cells.Item(135, 8).Value = 1348.5
$ws.Cells.Item(135, 9).Value = 1312.5
$ws.Cells.Item(135, 11).Value = 11812.5
$ws.Cells.Item(135, 13).Value = -9277.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 580.96155
$ws.Cells.Item(2, 9).Value = 580.8421
$ws.Cells.Item(2, 10).Value = 581.2857
$ws.Cells.Item(2, 11).Value = 580.8421
$ws.Cells.Item(2, 12).Value = 581.2857
$ws.Cells.Item(2, 13).Value = -467.8421
$ws.Cells.Item(2, 14).Value = -807.2857
$ws.Cells.Item(32, 8).Value = 3694.76
$ws.Cells.Item(32, 9).Value = 2120.95
$ws.Cells.Item(32, 10).Value = 9990
$ws.Cells.Item(32, 11).Value = 2120.95
$ws.Cells.Item(32, 12).Value = 9990
$ws.Cells.Item(32, 13).Value = -1833.95
$ws.Cells.Item(32, 14).Value = -10564
$ws.Cells.Item(61, 8).Value = 1733.9714
$ws.Cells.Item(61, 9).Value = 1157.3462
$ws.Cells.Item(61, 11).Value = 1157.3462
$ws.Cells.Item(61, 13).Value = -945.3462
$ws.Cells.Item(74, 8).Value = 3469.1316
$ws.Cells.Item(74, 9).Value = 3390.2903
$ws.Cells.Item(74, 11).Value = 3390.2903
$ws.Cells.Item(74, 13).Value = -2516.2903
$ws.Cells.Item(77, 8).Value = 3469.1316
$ws.Cells.Item(77, 9).Value = 3390.2903
$ws.Cells.Item(77, 11).Value = 16951.4515
$ws.Cells.Item(77, 13).Value = -12583.4515
$ws.Cells.Item(116, 8).Value = 580.96155
$ws.Cells.Item(116, 9).Value = 580.8421
$ws.Cells.Item(116, 10).Value = 581.2857
$ws.Cells.Item(116, 11).Value = 580.8421
$ws.Cells.Item(116, 12).Value = 581.2857
$ws.Cells.Item(116, 13).Value = 1713.1579
$ws.Cells.Item(116, 14).Value = -5169.2857
$ws.Cells.Item(122, 8).Value = 3090.7693
$ws.Cells.Item(122, 9).Value = 1528.5714
$ws.Cells.Item(122, 11).Value = 4585.7142
$ws.Cells.Item(122, 13).Value = -2135.7142
$ws.Cells.Item(136, 8).Value = 1733.9714
$ws.Cells.Item(136, 9).Value = 1157.3462
$ws.Cells.Item(136, 11).Value = 3472.0386
$ws.Cells.Item(136, 13).Value = -922.0385999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 580.96155
$ws.Cells.Item(3, 9).Value = 580.8421
$ws.Cells.Item(3, 10).Value = 581.2857
$ws.Cells.Item(3, 11).Value = 580.8421
$ws.Cells.Item(3, 12).Value = 581.2857
$ws.Cells.Item(3, 13).Value = -466.8421
$ws.Cells.Item(3, 14).Value = -809.2857
$ws.Cells.Item(99, 8).Value = 2638.125
$ws.Cells.Item(99, 9).Value = 1196.1904
$ws.Cells.Item(99, 11).Value = 1196.1904
$ws.Cells.Item(99, 13).Value = 301.8096
$ws.Cells.Item(105, 8).Value = 2554.9285
$ws.Cells.Item(105, 9).Value = 2537.925
$ws.Cells.Item(105, 11).Value = 2537.925
$ws.Cells.Item(105, 13).Value = -790.9250000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(63, 8).Value = 49995
$ws.Cells.Item(63, 10).Value = 49995
$ws.Cells.Item(63, 12).Value = 49995
$ws.Cells.Item(63, 14).Value = -51367
$ws.Cells.Item(66, 8).Value = 49995
$ws.Cells.Item(66, 10).Value = 49995
$ws.Cells.Item(66, 12).Value = 149985
$ws.Cells.Item(66, 14).Value = -156849
$ws.Cells.Item(132, 8).Value = 2768.1345
$ws.Cells.Item(132, 9).Value = 2057.7878
$ws.Cells.Item(132, 10).Value = 4001.8948
$ws.Cells.Item(132, 11).Value = 6173.3634
$ws.Cells.Item(132, 12).Value = 12005.6844
$ws.Cells.Item(132, 13).Value = -3643.3634
$ws.Cells.Item(132, 14).Value = -17065.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 663.7193
$ws.Cells.Item(113, 9).Value = 592.4146
$ws.Cells.Item(113, 10).Value = 846.4375
$ws.Cells.Item(113, 11).Value = 1777.2438
$ws.Cells.Item(113, 12).Value = 2539.3125
$ws.Cells.Item(113, 13).Value = 392.7562000000003
$ws.Cells.Item(113, 14).Value = -6879.3125
$ws.Cells.Item(131, 8).Value = 12196141
$ws.Cells.Item(131, 10).Value = 848.80554
$ws.Cells.Item(131, 12).Value = 2546.41662
$ws.Cells.Item(131, 14).Value = -12626.41662

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(105, 8).Value = 30000
$ws.Cells.Item(105, 10).Value = 30000
$ws.Cells.Item(105, 12).Value = 30000
$ws.Cells.Item(105, 14).Value = -36988
$ws.Cells.Item(132, 8).Value = 3233.625
$ws.Cells.Item(132, 9).Value = 1641.8334
$ws.Cells.Item(132, 11).Value = 4925.5002
$ws.Cells.Item(132, 13).Value = -2395.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4415.475
$ws.Cells.Item(40, 9).Value = 3985.6128
$ws.Cells.Item(40, 10).Value = 5896.1113
$ws.Cells.Item(40, 11).Value = 3985.6128
$ws.Cells.Item(40, 12).Value = 5896.1113
$ws.Cells.Item(40, 13).Value = -3849.6128
$ws.Cells.Item(40, 14).Value = -6168.1113
$ws.Cells.Item(46, 8).Value = 1759.1852
$ws.Cells.Item(46, 9).Value = 1588.1765
$ws.Cells.Item(46, 10).Value = 2049.9
$ws.Cells.Item(46, 11).Value = 1588.1765
$ws.Cells.Item(46, 12).Value = 2049.9
$ws.Cells.Item(46, 13).Value = -1400.1765
$ws.Cells.Item(46, 14).Value = -2425.9
$ws.Cells.Item(136, 8).Value = 4513.44
$ws.Cells.Item(136, 9).Value = 1671.6154
$ws.Cells.Item(136, 10).Value = 7592.0835
$ws.Cells.Item(136, 11).Value = 5014.8462
$ws.Cells.Item(136, 12).Value = 22776.2505
$ws.Cells.Item(136, 13).Value = -2464.8462
$ws.Cells.Item(136, 14).Value = -27876.2505
$ws.Cells.Item(140, 8).Value = 65345.91
$ws.Cells.Item(140, 10).Value = 65345.91
$ws.Cells.Item(140, 12).Value = 65345.91
$ws.Cells.Item(140, 14).Value = -75705.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 715.5238000000001
$ws.Cells.Item(107, 9).Value = 630.0625
$ws.Cells.Item(107, 10).Value = 989
$ws.Cells.Item(107, 11).Value = 1890.1875
$ws.Cells.Item(107, 12).Value = 2967
$ws.Cells.Item(107, 13).Value = 29.8125
$ws.Cells.Item(107, 14).Value = -6807
$ws.Cells.Item(113, 8).Value = 244
$ws.Cells.Item(113, 9).Value = 244
$ws.Cells.Item(113, 11).Value = 732
$ws.Cells.Item(113, 13).Value = 1438
$ws.Cells.Item(132, 8).Value = 6804264
$ws.Cells.Item(132, 9).Value = 622.0968
$ws.Cells.Item(132, 10).Value = 18521648
$ws.Cells.Item(132, 11).Value = 1866.2904
$ws.Cells.Item(132, 12).Value = 55564944
$ws.Cells.Item(132, 13).Value = 663.7095999999999
$ws.Cells.Item(132, 14).Value = -55570004
